# Fix syntax highlighting XML element order and add automatic text color
# contrast for shapes.
#
# For every auto-shape on the affected slides whose text still uses the
# original default run size (18pt / sz="1800"), this:
#   1. Turns on "Shrink text on overflow" (PpAutoSize.ppAutoSizeTextToFitShape)
#      so PowerPoint writes <a:normAutofit/> into <a:bodyPr>.
#   2. Computes a new font size that fits the shape's width/height using the
#      same width/height based fit heuristic the authoring tool uses:
#        width_based_pt  = 2.0 * shapeWidthPt  / longestLineLength
#        height_based_pt = 0.8 * shapeHeightPt / numberOfLines
#        fontPt = min(width_based_pt, height_based_pt, 44)  # 44pt hard cap
#      truncated to hundredths of a point (centipoints), matching PowerPoint's
#      sz unit.
#   3. Picks a contrasting run text color (white on dark fills, black on
#      light fills) using the sRGB luminance of the shape's own solid fill,
#      with a luminance threshold of 128.

function Get-Luminance([int]$r, [int]$g, [int]$b) {
    return (0.299 * $r) + (0.587 * $g) + (0.114 * $b)
}

function Test-IsDarkColor([int]$rgbInt) {
    $r = $rgbInt -band 0xFF
    $g = ($rgbInt -shr 8) -band 0xFF
    $b = ($rgbInt -shr 16) -band 0xFF
    $lum = Get-Luminance $r $g $b
    return $lum -lt 128
}

function Get-TextColorRgb([int]$fillRgbInt) {
    if (Test-IsDarkColor $fillRgbInt) {
        # white
        return 255 + (255 * 256) + (255 * 65536)
    } else {
        # black
        return 0
    }
}

function Get-FitFontSizePt([double]$widthPt, [double]$heightPt, [string]$text) {
    $lines = $text -split "`n"
    $lineCount = $lines.Count
    $maxLen = 0
    foreach ($line in $lines) {
        if ($line.Length -gt $maxLen) { $maxLen = $line.Length }
    }
    if ($maxLen -lt 1) { $maxLen = 1 }

    $widthBased = (2.0 * $widthPt) / $maxLen
    $heightBased = (0.8 * $heightPt) / $lineCount

    $fontPt = [Math]::Min($widthBased, $heightBased)
    $fontPt = [Math]::Min($fontPt, 44)

    # Truncate (floor) to hundredths of a point, like sz centipoints.
    $centi = [Math]::Floor($fontPt * 100)
    return ($centi / 100.0)
}

$p = $ppt.ActivePresentation

# Walk every slide/shape in the deck; only act on plain auto-shapes
# ("Shape N") that still carry the original, un-autofitted 18pt runs and a
# solid background fill - i.e. exactly the shapes the generator originally
# emitted via `generate_text_xml_with...` without color-contrast support.
for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }
        $tf = $sh.TextFrame
        if (-not $tf.HasText) { continue }
        if ($sh.Name -notlike "Shape *") { continue }

        $tr = $tf.TextRange
        if ($tr.Font.Size -ne 18) { continue }
        if (-not $sh.Fill.Visible) { continue }

        $fillRgb = $sh.Fill.ForeColor.RGB

        $newSizePt = Get-FitFontSizePt $sh.Width $sh.Height $tr.Text
        $newColorRgb = Get-TextColorRgb $fillRgb

        # Shrink text on overflow -> <a:normAutofit/>
        $tf.AutoSize = 2

        $tr.Font.Size = $newSizePt
        $tr.Font.Color.RGB = $newColorRgb
    }
}
